# Se cambia la extension a jpg
# The sheet has a grid of CONCATENATE(...) formulas in G1:K22 that build
# image-file-name strings like '0B.png', from the values in columns A:E.
# The edit simply swaps the ".png" extension for ".jpg" everywhere those
# formulas (and their cached string results) reference it.
#
# Using Range.Replace against the formula text (rather than rewriting each
# cell's .Formula individually) mirrors doing a Find & Replace across the
# sheet in Excel: it rewrites only the formula source text in place and
# leaves every shared-formula grouping, cell type and untouched cell alone
# (no phantom cells get created in rows where a column was previously
# empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("G1:K22")
$rng.Replace(".png", ".jpg") | Out-Null
